$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4963948798677646
$ws.Range("C2").Value = 0.9939133157012632
$ws.Range("D2").Value = 0.3985919732244694
$ws.Range("E2").Value = 0.7627277865358761
$ws.Range("F2").Value = 1.184322899759067
$ws.Range("G2").Value = 1.034839871199764
$ws.Range("H2").Value = 1.088646979932846
$ws.Range("B3").Value = 0.4975184358334985
$ws.Range("C3").Value = -0.09780290664329527
$ws.Range("D3").Value = 0.2663329066681115
$ws.Range("E3").Value = 0.6879280198913023
$ws.Range("F3").Value = 0.5384449913319997
$ws.Range("G3").Value = 0.5922521000650818
$ws.Range("B4").Value = -0.5953213424767938
$ws.Range("C4").Value = -0.231185529165387
$ws.Range("D4").Value = 0.1904095840578037
$ws.Range("E4").Value = 0.04092655549850122
$ws.Range("F4").Value = 0.09473366423158319
$ws.Range("G4").Value = -0.2300153708986841
$ws.Range("H4").Value = 0.1053516340798524
$ws.Range("I4").Value = 0.1988927554601076
$ws.Range("J4").Value = -0.2839541438535775
$ws.Range("B5").Value = 0.3641358133114068
$ws.Range("C5").Value = 0.7857309265345975
$ws.Range("D5").Value = 0.636247897975295
$ws.Range("E5").Value = 0.6900550067083771
$ws.Range("F5").Value = 0.3653059715781097
$ws.Range("G5").Value = 0.7006729765566462
$ws.Range("H5").Value = 0.7942140979369015
$ws.Range("I5").Value = 0.3113671986232163
$ws.Range("B6").Value = 0.4215951132231908
$ws.Range("C6").Value = 0.2721120846638883
$ws.Range("D6").Value = 0.3259191933969702
$ws.Range("E6").Value = 0.00117015826670297
$ws.Range("F6").Value = 0.3365371632452395
$ws.Range("G6").Value = 0.4300782846254947
$ws.Range("H6").Value = -0.0527686146881905
$ws.Range("B7").Value = -0.1494830285593025
$ws.Range("C7").Value = -0.09567591982622053
$ws.Range("D7").Value = -0.4204249549564878
$ws.Range("E7").Value = -0.0850579499779513
$ws.Range("F7").Value = 0.008483171402303896
$ws.Range("G7").Value = -0.4743637279113813
$ws.Range("B8").Value = 0.05380710873308198
$ws.Range("C8").Value = -0.2709419263971853
$ws.Range("D8").Value = 0.06442507858135121
$ws.Range("E8").Value = 0.1579661999616064
$ws.Range("F8").Value = -0.3248806993520788
$ws.Range("G8").Value = -0.2152921167545969
$ws.Range("H8").Value = -0.186132674248455
$ws.Range("I8").Value = -0.3848271424380556
$ws.Range("B9").Value = -0.3247490351302673
$ws.Range("C9").Value = 0.01061796984826924
$ws.Range("D9").Value = 0.1041590912285244
$ws.Range("E9").Value = -0.3786878080851607
$ws.Range("F9").Value = -0.2690992254876789
$ws.Range("G9").Value = -0.239939782981537
$ws.Range("H9").Value = -0.4386342511711376
$ws.Range("B10").Value = 0.3353670049785365
$ws.Range("C10").Value = 0.4289081263587917
$ws.Range("D10").Value = -0.05393877295489347
$ws.Range("E10").Value = 0.05564980964258837
$ws.Range("F10").Value = 0.08480925214873025
$ws.Range("G10").Value = -0.1138852160408703
$ws.Range("B11").Value = 0.09354112138025519
$ws.Range("C11").Value = -0.38930577793343
$ws.Range("D11").Value = -0.2797171953359481
$ws.Range("E11").Value = -0.2505577528298062
$ws.Range("F11").Value = -0.4492522210194068
$ws.Range("B12").Value = -0.4828468993136852
$ws.Range("C12").Value = -0.3732583167162034
$ws.Range("D12").Value = -0.3440988742100615
$ws.Range("E12").Value = -0.542793342399662
$ws.Range("B13").Value = 0.1095885825974818
$ws.Range("C13").Value = 0.1387480251036237
$ws.Range("D13").Value = -0.05994644308597685
$ws.Range("B14").Value = 0.02915944250614189
$ws.Range("C14").Value = -0.1695350256834587
$ws.Range("B15").Value = -0.1986944681896006
